$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.216.56"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "3.065.52"
$ws.Range("E3").Value = "  +2.77%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'527.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.92%  "
$ws.Range("D6").Value = "'143.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +5.50%  "
$ws.Range("D9").Value = "'7.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.18%  "
$ws.Range("E10").Value = "  +7.51%  "
$ws.Range("E11").Value = "  +6.14%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").Value = "3.592.03"
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("D14").Value = "'27.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.25%  "
$ws.Range("D15").Value = "'0.0000174"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +16.18%  "
$ws.Range("D16").Value = "58.143.26"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "'6.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +10.20%  "
$ws.Range("D18").Value = "3.070.11"
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").Value = "'13.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.26%  "
$ws.Range("E20").Value = "  +5.52%  "
$ws.Range("D21").Value = "'343.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.26%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'0.507"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.35%  "
$ws.Range("D25").Value = "'65.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.67%  "
$ws.Range("D26").Value = "0.0₃0974"
$ws.Range("E26").Value = "  +9.97%  "
$ws.Range("E27").Value = "  +3.84%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +9.92%  "
$ws.Range("D30").Value = "'7.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.85%  "
$ws.Range("E31").Value = "  +7.38%  "
$ws.Range("E32").Value = "  +6.02%  "
$ws.Range("D33").Value = "'21.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.71%  "
$ws.Range("E34").Value = "  +9.09%  "
$ws.Range("D35").Value = "'158.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").Value = "'5.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.96%  "
$ws.Range("E37").Value = "  +4.01%  "
$ws.Range("D38").Value = "'26.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.70%  "
$ws.Range("D39").Value = "'0.0701"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.38%  "
$ws.Range("D40").Value = "3.100.75"
$ws.Range("E40").Value = "  +2.78%  "
$ws.Range("D41").Value = "'37.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'3.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.28%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.668"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +6.10%  "
$ws.Range("D46").Value = "2.343.41"
$ws.Range("E46").Value = "  +5.18%  "
$ws.Range("E47").Value = "  +4.75%  "
$ws.Range("E48").Value = "  +4.64%  "
$ws.Range("D49").Value = "'6.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.54%  "
$ws.Range("D50").Value = "'0.0244"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.82%  "
$ws.Range("E51").Value = "  +6.81%  "
